$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sreyas")

# New match rows appended to the bottom of the table (rows 39-47)
$rows = @(
    @{ Row = 39; A = "Sparta Rotterdam U15";        B = "Ajax JO15-1";            C = 2037; D = 90;              E = 45411; F = 500 },
    @{ Row = 40; A = "Nita FA";                      B = "Pune Krida PA";          C = 2000; D = 90;              E = 45413; F = 500 },
    @{ Row = 41; A = "Roots Red U15";                B = "Dash Athlectics U15";    C = 2036; D = "Less than 60";  E = 45414; F = 300 },
    @{ Row = 42; A = "Young Blues Elite FC U17";      B = "Raman SA U17";           C = 2072; D = "Less than 60";  E = 45416; F = 300 },
    @{ Row = 43; A = "Ajax JO15-1";                  B = "Academy Fukushima U15";  C = 2074; D = 90;              E = 45417; F = 500 },
    @{ Row = 44; A = "Loco. Tbilisi";                B = "FC Gareji Sagarejo";     C = 2089; D = 90;              E = 45420; F = 500 },
    @{ Row = 45; A = "RFYC Development Squad";       B = "Mumbai City FC B";       C = 1991; D = "Less than 60";  E = 45421; F = 300 },
    @{ Row = 46; A = "FC Utrecht U17";               B = "Ajax JO17-1";            C = 2097; D = 90;              E = 45422; F = 500 },
    @{ Row = 47; A = "Feyenoord O17";                B = " Aax JO17-1";            C = 2107; D = 90;              E = 45425; F = 500 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 5).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($rowNum, 6).Value = $r.F
}
